$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 5.586269137925634)
    3 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 5.586269137925634)
    4 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 0, 8.974608811992548)
    5 = @(0.2917716402565462, 0.306821227259698, 22.3905356188092, 10.19245300693656, 1, 33.181581493262)
    6 = @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 1, 13.26162679800665)
    7 = @(0.2917716402565462, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 0, 12.89274599722724)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}
